$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DistanceMatrix")

# The surviving rows are the original rows 1, 3 and 5 (header + the two
# "car / en-EN" scenario rows). Delete the others, working bottom-up so
# row indices of rows not yet processed stay valid.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(2).Delete()

# Delete column J (the "OS" column)
$ws.Columns.Item(10).Delete()

# Update selection to match the after-state (F9)
$ws.Range("F9").Select()
